# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F2").Value = 14955
$sheet1.Range("F3").Value = 18813
$sheet1.Range("F13").Value = 55
$sheet1.Range("F14").Value = 132
$sheet1.Range("F17").Value = 1438
$sheet1.Range("F22").Value = 7819
$sheet1.Range("F25").Value = 59
$sheet1.Range("F26").Value = 1234
$sheet1.Range("F28").Value = 6011
$sheet1.Range("F31").Value = 166
$sheet1.Range("F34").Value = 5377

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F2").Value = 14955
$sheet4.Range("F3").Value = 18813
$sheet4.Range("F13").Value = 55
$sheet4.Range("F14").Value = 132
$sheet4.Range("F17").Value = 1438
$sheet4.Range("F23").Value = 7819
$sheet4.Range("F26").Value = 59
$sheet4.Range("F27").Value = 1234
$sheet4.Range("F31").Value = 6011
$sheet4.Range("F34").Value = 166
$sheet4.Range("F37").Value = 5377
